# "Generate Report for Archive"
# 1) Status text updated from "Ready for handoff" -> "In Translation"
#    on the Overview sheet (zh-cn / de-de status columns, E & F) and on
#    the per-locale "zh-cn" / "de-de" sheets (Status column, C).
# 2) The Status column(s) got narrower to fit the new (shorter) text:
#    width 17.2159881591797 -> 13.4101848602295 (Overview!E:F, zh-cn!C, de-de!C)

$wb = $excel.ActiveWorkbook

$oldStatus = "Ready for handoff"
$newStatus = "In Translation"

# --- Overview sheet: zh-cn / de-de status cells (columns E & F) ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Range("E3").Value = $newStatus
$wsOverview.Range("F3").Value = $newStatus

# --- zh-cn sheet: Status column (C) ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Range("C3").Value = $newStatus

# --- de-de sheet: Status column (C) ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Range("C3").Value = $newStatus

# --- Shrink the now-narrower Status columns to match the new text ---
# ColumnWidth is in characters; 12.5 is the closest input that lands the
# saved (pixel-rounded) column width nearest the target 13.4101848602295.
$newColumnWidth = 12.5
$wsOverview.Range("E1").ColumnWidth = $newColumnWidth
$wsOverview.Range("F1").ColumnWidth = $newColumnWidth
$wsZhCn.Range("C1").ColumnWidth = $newColumnWidth
$wsDeDe.Range("C1").ColumnWidth = $newColumnWidth
